{"js": "// Office.js (Word JavaScript API) edit script\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// ---- Change 1: consolidate the Synopsis summary paragraph's runs into one run ----\n// (text content is unchanged; only the run structure is simplified)\nconst synopsisResults = body.search(\"The Borough areas of London, near the bridge\", { matchCase: true, matchWholeWord: false });\nsynopsisResults.load(\"items\");\nawait context.sync();\n\nif (synopsisResults.items.length === 0) {\n  throw new Error(\"Could not find the Synopsis summary paragraph\");\n}\n\nconst synopsisParagraph = synopsisResults.items[0].paragraphs.getFirst();\nconst newSynopsisText = \"The Borough areas of London, near the bridge. George is with Soula, Elektra, Melody, Fleur, Susan and Michael at a pavement cafe. A truck starts driving into pedestrians. George, Soula and Elektra stop it and capture the terrorists.\";\nsynopsisParagraph.insertText(newSynopsisText, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---- Change 2: expand the \"Car crash...\" paragraph with the full character\n// descriptions, then add a new paragraph holding the old \"{Description of each}\"\n// placeholder text (replacing the blank paragraph that used to directly follow it) ----\nconst carCrashResults = body.search(\"Car crash sounds and screams\", { matchCase: true, matchWholeWord: false });\ncarCrashResults.load(\"items\");\nawait context.sync();\n\nif (carCrashResults.items.length === 0) {\n  throw new Error(\"Could not find the Car crash paragraph\");\n}\n\nconst carCrashParagraph = carCrashResults.items[0].paragraphs.getFirst();\nconst newCarCrashText = \"Car crash sounds and screams. Too distant for ordinary humans to hear. One table at an outdoor cafe, every ear pricked up. The first was currently known as George Kominos, Greek looking, early thirties, black hair that fell nearly to his shoulders, three day stubble, dressed all in black, including a black broad brimmed leather hat, and mirror sunshades. Three hundred year old master vampire and day walker. Next to him was Soula, currently Soula Kominos, flame red hair caught in a braid that descended past her shoulders, pale skin, blue eyes, clad in green slacks and matching top. Looked about twenty five, except she was George\u2019s number one human servant, possessed of similar powers to George himself, and older than him. On George\u2019s other side was Elektra, currently Elektra Kominos, a blue eyed blonde who looked about twenty. Her long hair fell in waves past her shoulders, held in place with jewelled combs, and she wore medium blue slacks and matching top. She was George\u2019s number two human servant, perhaps fifteen years younger than George. Next to Soula sat Melody Waters, milk pale skin, raven black hair, leaf green eyes. She wore brown slacks and a leaf green top, there was silver wire bound into her locks, dividing her tresses into nine braids of shoulder length. Her ears were pointy, denoting her pure Fae blood, though the braids hid most of this. She appeared  to be in her mid twenties, though this was quite deceiving. Next to Elektra sat Fleur Kominos, seemingly seventeen, with shoulder length auburn hair, milk pale skin and sky blue eyes. She seemed mature for her age, but quite human, though that was deceiving. No one at the table was sure what she was really, except that George had found her as a homeless young child thirty years ago who was producing flames from her finger tips, and invited her home. His women had adopted her, and subsequently they\u2019d all met Melody when they sent Fleur to school. Fleur was wearing a white top and dark blue slacks. Next to Fleur was Michael, and next to him was Susan, his partner. Michael was slim with brown hair and eyes, dressed in a faded denim shirt, faded jeans and sneakers. Susan was blonde, hair cut in a bob, with blue eyes, pale skin with heavy makeup, a blue top that was the same design as Fleur\u2019s and the same colour as Elektra\u2019s, designer dark blue jeans and sneakers. All the others were wearing comfortable walking shoes.\";\ncarCrashParagraph.insertText(newCarCrashText, Word.InsertLocation.replace);\nawait context.sync();\n\n// The paragraph that used to immediately follow the \"Car crash...\" paragraph was\n// blank; remove it and insert a new paragraph with the placeholder text instead.\nconst followingParagraph = carCrashParagraph.getNextOrNullObject();\nawait context.sync();\nfollowingParagraph.load(\"text\");\nawait context.sync();\n\nconst placeholderText = \"{Description of each} \";\nif (!followingParagraph.isNullObject && followingParagraph.text === \"\") {\n  followingParagraph.insertText(placeholderText, Word.InsertLocation.replace);\n} else {\n  carCrashParagraph.insertParagraph(placeholderText, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script\n# $word.ActiveDocument is available as $d\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexForRange($doc, $targetRange) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($p.Range.Start -le $targetRange.Start -and $p.Range.End -ge $targetRange.End) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# ---- Change 1: consolidate the Synopsis summary paragraph's runs into one run ----\n# (text content is unchanged; only the run structure is simplified)\n$findRange1 = $d.Content\n$findRange1.Find.ClearFormatting()\n$found1 = $findRange1.Find.Execute(\"The Borough areas of London, near the bridge\")\nif (-not $found1) {\n    throw \"Could not find the Synopsis summary paragraph\"\n}\n$synopsisIdx = Get-ParagraphIndexForRange $d $findRange1\n$synopsisRange = $d.Paragraphs.Item($synopsisIdx).Range\n$synopsisRange.MoveEnd(1, -1) | Out-Null\n$synopsisRange.Text = \"The Borough areas of London, near the bridge. George is with Soula, Elektra, Melody, Fleur, Susan and Michael at a pavement cafe. A truck starts driving into pedestrians. George, Soula and Elektra stop it and capture the terrorists.\"\n\n# ---- Change 2: expand the \"Car crash...\" paragraph with the full character\n# descriptions, then fill the blank paragraph that directly follows it with the\n# old \"{Description of each}\" placeholder text ----\n$findRange2 = $d.Content\n$findRange2.Find.ClearFormatting()\n$found2 = $findRange2.Find.Execute(\"Car crash sounds and screams\")\nif (-not $found2) {\n    throw \"Could not find the Car crash paragraph\"\n}\n$carCrashIdx = Get-ParagraphIndexForRange $d $findRange2\n$carCrashRange = $d.Paragraphs.Item($carCrashIdx).Range\n$carCrashRange.MoveEnd(1, -1) | Out-Null\n$carCrashRange.Text = \"Car crash sounds and screams. Too distant for ordinary humans to hear. One table at an outdoor cafe, every ear pricked up. The first was currently known as George Kominos, Greek looking, early thirties, black hair that fell nearly to his shoulders, three day stubble, dressed all in black, including a black broad brimmed leather hat, and mirror sunshades. Three hundred year old master vampire and day walker. Next to him was Soula, currently Soula Kominos, flame red hair caught in a braid that descended past her shoulders, pale skin, blue eyes, clad in green slacks and matching top. Looked about twenty five, except she was George\u2019s number one human servant, possessed of similar powers to George himself, and older than him. On George\u2019s other side was Elektra, currently Elektra Kominos, a blue eyed blonde who looked about twenty. Her long hair fell in waves past her shoulders, held in place with jewelled combs, and she wore medium blue slacks and matching top. She was George\u2019s number two human servant, perhaps fifteen years younger than George. Next to Soula sat Melody Waters, milk pale skin, raven black hair, leaf green eyes. She wore brown slacks and a leaf green top, there was silver wire bound into her locks, dividing her tresses into nine braids of shoulder length. Her ears were pointy, denoting her pure Fae blood, though the braids hid most of this. She appeared  to be in her mid twenties, though this was quite deceiving. Next to Elektra sat Fleur Kominos, seemingly seventeen, with shoulder length auburn hair, milk pale skin and sky blue eyes. She seemed mature for her age, but quite human, though that was deceiving. No one at the table was sure what she was really, except that George had found her as a homeless young child thirty years ago who was producing flames from her finger tips, and invited her home. His women had adopted her, and subsequently they\u2019d all met Melody when they sent Fleur to school. Fleur was wearing a white top and dark blue slacks. Next to Fleur was Michael, and next to him was Susan, his partner. Michael was slim with brown hair and eyes, dressed in a faded denim shirt, faded jeans and sneakers. Susan was blonde, hair cut in a bob, with blue eyes, pale skin with heavy makeup, a blue top that was the same design as Fleur\u2019s and the same colour as Elektra\u2019s, designer dark blue jeans and sneakers. All the others were wearing comfortable walking shoes.\"\n\n$followingPara = $d.Paragraphs.Item($carCrashIdx + 1)\n$followingRange = $followingPara.Range\n$isBlankParagraph = (($followingRange.End - $followingRange.Start) -eq 1)\nif ($isBlankParagraph) {\n    $followingRange.MoveEnd(1, -1) | Out-Null\n    $followingRange.Text = \"{Description of each} \"\n} else {\n    $carCrashRange2 = $d.Paragraphs.Item($carCrashIdx).Range\n    $carCrashRange2.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($carCrashIdx + 1)\n    $newRange = $newPara.Range\n    $newRange.MoveEnd(1, -1) | Out-Null\n    $newRange.Text = \"{Description of each} \"\n}\n"}
